$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.882.54"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").Value = "2.580.36"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.82%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").Value = "2.585.70"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.114"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").Value = "3.047.86"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("D17").Value = "62.783.80"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").Value = "2.598.09"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "334.48"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.23%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000109"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "564.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.98"
$ws.Range("D33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.68"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.39"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.22"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.396"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.38"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "155.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.38"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0576"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.66%  "

$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0986"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("E50").Value = "  -2.26%  "

$ws.Range("D51").Value = "0.0₆0228"
$ws.Range("E51").Value = "  -1.52%  "
